$wb = $excel.ActiveWorkbook

# Delete the two "(2)" duplicate sheets
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("5-state confusion matrix (2)").Delete()
$wb.Worksheets.Item("2-state confusion matrix (2)").Delete()

# Rename "Tabelle1" to "Results Table"
$wb.Worksheets.Item("Tabelle1").Name = "Results Table"

# Make "Results Table" the active sheet/tab and set its selection
$resultsSheet = $wb.Worksheets.Item("Results Table")
$resultsSheet.Activate()
$resultsSheet.Range("G27").Select()
